$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update probability matrix values per games pulled March 7

# Row 2
$ws.Range("B2").Value = 0.2162162162162162
$ws.Range("C2").Value = 0.5297297297297298
$ws.Range("J2").Value = 0.02432432432432433
$ws.Range("P2").Value = 0.1486486486486487
$ws.Range("S2").Value = 0.08108108108108109

# Row 3
$ws.Range("B3").Value = 0.01932367149758454
$ws.Range("C3").Value = 0.05797101449275362
$ws.Range("J3").Value = 0.02898550724637681
$ws.Range("P3").Value = 0.7294685990338164
$ws.Range("S3").Value = 0.1642512077294686

# Row 4
$ws.Range("J4").Value = 0.01724137931034483
$ws.Range("P4").Value = 0.8103448275862069
$ws.Range("S4").Value = 0.1724137931034483

# Row 5
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333

# Row 6
$ws.Range("B6").Value = 0.05809128630705394
$ws.Range("D6").Value = 0.008298755186721992
$ws.Range("E6").Value = 0.004149377593360996
$ws.Range("F6").Value = 0.07883817427385892
$ws.Range("J6").Value = 0.2697095435684647
$ws.Range("O6").Value = 0.03319502074688797
$ws.Range("Q6").Value = 0.1244813278008299
$ws.Range("R6").Value = 0.07053941908713693
$ws.Range("S6").Value = 0.3526970954356847

# Row 7
$ws.Range("B7").Value = 0.08518518518518518
$ws.Range("D7").Value = 0.03333333333333333
$ws.Range("F7").Value = 0.03703703703703703
$ws.Range("J7").Value = 0.1185185185185185
$ws.Range("O7").Value = 0.04444444444444445
$ws.Range("Q7").Value = 0.137037037037037
$ws.Range("R7").Value = 0.1222222222222222
$ws.Range("S7").Value = 0.4222222222222222

# Row 8
$ws.Range("B8").Value = 0.09724770642201835
$ws.Range("D8").Value = 0.02568807339449541
$ws.Range("E8").Value = 0.001834862385321101
$ws.Range("F8").Value = 0.05321100917431193
$ws.Range("J8").Value = 0.1302752293577982
$ws.Range("O8").Value = 0.01834862385321101
$ws.Range("Q8").Value = 0.1357798165137615
$ws.Range("R8").Value = 0.07155963302752294
$ws.Range("S8").Value = 0.4660550458715597

# Row 9
$ws.Range("B9").Value = 0.09174311926605505
$ws.Range("D9").Value = 0.01376146788990826
$ws.Range("F9").Value = 0.06880733944954129
$ws.Range("J9").Value = 0.1055045871559633
$ws.Range("O9").Value = 0.01834862385321101
$ws.Range("Q9").Value = 0.1559633027522936
$ws.Range("R9").Value = 0.09174311926605505
$ws.Range("S9").Value = 0.4541284403669725

# Row 10
$ws.Range("B10").Value = 0.1170798898071625
$ws.Range("D10").Value = 0.02134986225895317
$ws.Range("E10").Value = 0.0006887052341597796
$ws.Range("F10").Value = 0.06473829201101929
$ws.Range("J10").Value = 0.1129476584022039
$ws.Range("O10").Value = 0.02341597796143251
$ws.Range("Q10").Value = 0.196969696969697
$ws.Range("R10").Value = 0.08264462809917356
$ws.Range("S10").Value = 0.3801652892561984

# Row 11
$ws.Range("G11").Value = 0.1445221445221445
$ws.Range("J11").Value = 0.0979020979020979
$ws.Range("K11").Value = 0.2004662004662005
$ws.Range("L11").Value = 0.5524475524475524
$ws.Range("S11").Value = 0.004662004662004662

# Row 12
$ws.Range("G12").Value = 0.7302904564315352
$ws.Range("J12").Value = 0.2323651452282158
$ws.Range("K12").Value = 0.004149377593360996
$ws.Range("L12").Value = 0.008298755186721992
$ws.Range("S12").Value = 0.02489626556016597

# Row 13
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.2115384615384615
$ws.Range("S13").Value = 0.07692307692307693

# Row 15
$ws.Range("F15").Value = 0.03448275862068965
$ws.Range("H15").Value = 0.1637931034482759
$ws.Range("I15").Value = 0.0603448275862069
$ws.Range("J15").Value = 0.3146551724137931
$ws.Range("K15").Value = 0.05603448275862069
$ws.Range("M15").Value = 0.01293103448275862
$ws.Range("N15").Value = 0.004310344827586207
$ws.Range("O15").Value = 0.04310344827586207
$ws.Range("S15").Value = 0.3103448275862069

# Row 16
$ws.Range("F16").Value = 0.01606425702811245
$ws.Range("H16").Value = 0.1887550200803213
$ws.Range("I16").Value = 0.06425702811244979
$ws.Range("J16").Value = 0.4216867469879518
$ws.Range("K16").Value = 0.144578313253012
$ws.Range("M16").Value = 0.004016064257028112
$ws.Range("O16").Value = 0.02811244979919679
$ws.Range("S16").Value = 0.1325301204819277

# Row 17
$ws.Range("F17").Value = 0.02173913043478261
$ws.Range("H17").Value = 0.2173913043478261
$ws.Range("I17").Value = 0.08913043478260869
$ws.Range("J17").Value = 0.3717391304347826
$ws.Range("K17").Value = 0.1217391304347826
$ws.Range("M17").Value = 0.01304347826086956
$ws.Range("N17").Value = 0.006521739130434782
$ws.Range("O17").Value = 0.04347826086956522
$ws.Range("S17").Value = 0.1152173913043478

# Row 18
$ws.Range("F18").Value = 0.01310043668122271
$ws.Range("H18").Value = 0.1965065502183406
$ws.Range("I18").Value = 0.08733624454148471
$ws.Range("J18").Value = 0.3624454148471616
$ws.Range("K18").Value = 0.1397379912663755
$ws.Range("M18").Value = 0.02183406113537118
$ws.Range("O18").Value = 0.05676855895196507
$ws.Range("S18").Value = 0.1222707423580786

# Row 19
$ws.Range("F19").Value = 0.01291155584247902
$ws.Range("H19").Value = 0.2098127824402841
$ws.Range("I19").Value = 0.08327953518398967
$ws.Range("J19").Value = 0.3615235635894125
$ws.Range("K19").Value = 0.1304067140090381
$ws.Range("M19").Value = 0.02388637830858618
$ws.Range("N19").Value = 0.0006455777921239509
$ws.Range("O19").Value = 0.05551969012265978
$ws.Range("S19").Value = 0.1220142027114267
